# Update cryptos list with latest price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.153.56'
$ws.Range("E2").Value = '  +8.48%  '

# Row 3
$ws.Range("D3").Value = '3.520.66'
$ws.Range("E3").Value = '  +12.07%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '190.79'
$ws.Range("E5").Value = '  +13.57%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '550.09'
$ws.Range("E6").Value = '  +6.77%  '

# Row 7
$ws.Range("D7").Value = '3.511.36'
$ws.Range("E7").Value = '  +11.77%  '

# Row 8
$ws.Range("E8").Value = '  +3.91%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.04%  '

# Row 11
$ws.Range("E11").Value = '  +19.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.84'
$ws.Range("E12").Value = '  +6.43%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  +10.01%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.36'
$ws.Range("E14").Value = '  +6.01%  '

# Row 15
$ws.Range("D15").Value = '4.076.54'
$ws.Range("E15").Value = '  +12.07%  '

# Row 16
$ws.Range("D16").Value = '3.513.87'
$ws.Range("E16").Value = '  +12.21%  '

# Row 17
$ws.Range("E17").Value = '  +6.03%  '

# Row 18
$ws.Range("D18").Value = '67.142.90'
$ws.Range("E18").Value = '  +8.78%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.18'
$ws.Range("E19").Value = '  +7.96%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.86'
$ws.Range("E20").Value = '  +9.77%  '

# Row 21
$ws.Range("E21").Value = '  +4.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '426.91'
$ws.Range("E22").Value = '  +19.24%  '

# Row 23
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.92'
$ws.Range("E23").Value = '  +7.09%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.68'
$ws.Range("E24").Value = '  +6.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.18'
$ws.Range("E25").Value = '  +7.82%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.14'
$ws.Range("E26").Value = '  +0.79%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.90'
$ws.Range("E27").Value = '  +13.47%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.96'
$ws.Range("E28").Value = '  +8.74%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.88'
$ws.Range("E29").Value = '  +11.59%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.17'
$ws.Range("E30").Value = '  +8.92%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '653.81'
$ws.Range("E31").Value = '  +2.69%  '

# Row 32
$ws.Range("E32").Value = '  +5.67%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.71'
$ws.Range("E33").Value = '  +5.47%  '

# Row 34
$ws.Range("E34").Value = '  +8.08%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.43'
$ws.Range("E35").Value = '  +5.58%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.52'
$ws.Range("E36").Value = '  +7.49%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0816'
$ws.Range("E37").Value = '  +20.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.01%  '

# Row 39
$ws.Range("E39").Value = '  +6.21%  '

# Row 40
$ws.Range("E40").Value = '  +14.60%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.30'
$ws.Range("E41").Value = '  +16.06%  '

# Row 42
$ws.Range("E42").Value = '  +0.07%  '

# Row 43
$ws.Range("D43").Value = '3.011.49'
$ws.Range("E43").Value = '  +5.58%  '

# Row 44
$ws.Range("E44").Value = '  +5.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  +15.50%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.35'
$ws.Range("E46").Value = '  +12.59%  '

# Row 47
$ws.Range("E47").Value = '  +9.09%  '

# Row 48
$ws.Range("E48").Value = '  +4.01%  '

# Row 49
$ws.Range("E49").Value = '  +7.70%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.72'
$ws.Range("E50").Value = '  +17.69%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '140.29'
$ws.Range("E51").Value = '  +6.81%  '
